$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml)
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value = 205
$ws1.Range("F3").Value = 5526
$ws1.Range("F6").Value = 32
$ws1.Range("F7").Value = 654
$ws1.Range("F8").Value = 639
$ws1.Range("F9").Value = 5
$ws1.Range("F10").Value = 1075
$ws1.Range("F12").Value = 1532
$ws1.Range("F13").Value = 5062
$ws1.Range("F14").Value = 453
$ws1.Range("F15").Value = 236
$ws1.Range("F16").Value = 205
$ws1.Range("F17").Value = 14
$ws1.Range("F18").Value = 3
$ws1.Range("F19").Value = 108
$ws1.Range("F20").Value = 4334
$ws1.Range("F21").Value = 204
$ws1.Range("F22").Value = 1150
$ws1.Range("F23").Value = 119
$ws1.Range("F24").Value = 58
$ws1.Range("F26").Value = 55
$ws1.Range("F27").Value = 163
$ws1.Range("F28").Value = 61
$ws1.Range("F31").Value = 341
$ws1.Range("C32").Value = "南昌·ACG CLUB动漫游戏嘉年华"
$ws1.Range("F32").Value = 5

# Sheet "全部类型" (sheet4.xml)
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value = 205
$ws4.Range("F4").Value = 5526
$ws4.Range("F7").Value = 32
$ws4.Range("F8").Value = 654
$ws4.Range("F9").Value = 639
$ws4.Range("F10").Value = 5
$ws4.Range("F11").Value = 1075
$ws4.Range("F13").Value = 1532
$ws4.Range("F14").Value = 5062
$ws4.Range("F15").Value = 453
$ws4.Range("F16").Value = 236
$ws4.Range("F17").Value = 205
$ws4.Range("F18").Value = 14
$ws4.Range("F19").Value = 3
$ws4.Range("F20").Value = 108
$ws4.Range("F21").Value = 4334
$ws4.Range("F22").Value = 204
$ws4.Range("F23").Value = 1150
$ws4.Range("F24").Value = 119
$ws4.Range("F25").Value = 58
$ws4.Range("F27").Value = 55
$ws4.Range("F28").Value = 163
$ws4.Range("F29").Value = 61
$ws4.Range("F32").Value = 341
$ws4.Range("C33").Value = "南昌·ACG CLUB动漫游戏嘉年华"
$ws4.Range("F33").Value = 5
